# "Added View Event code changes"
# Update the AutoCreate/CreateAuto event rows on the EventDetails2 sheet and
# move the active selection to reflect where the user was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EventDetails2")
$ws.Activate()

# Row 2: event title + trial date
$ws.Range("A2").Value = "CreateAutoEvent"
$ws.Range("B2").Value = "10/29/2025"

# Row 3: event title + trial date (write B3 before A3 so new shared-string
# entries land in the same order as the authored workbook)
$ws.Range("B3").Value = "10/30/2025"
$ws.Range("A3").Value = "CreateAuto Event2"

# Move the selection/active cell like the author left it
$ws.Range("B10").Select()
